# "Generate Report for Handback"
# The handback report generator has run again: the 982bd13f-... file (and
# cc4d4c40-..., which depends on it) has now been handed back for both the
# zh-cn and de-de locales. This updates:
#   - Overview sheet: status "Ready for handoff" -> "Handed back: in sync with en-US"
#     for the two affected rows.
#   - zh-cn / de-de sheets: status updated the same way, plus the
#     "Latest Target File" / "Latest Handback File" columns (previously
#     empty) now get filled in, and "Latest Handback DateTime" moves from
#     the epoch placeholder to a real timestamp.

$wb = $excel.ActiveWorkbook

$HANDED_BACK = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = $HANDED_BACK
$ov.Range("C3").Value = $HANDED_BACK
$ov.Range("B4").Value = $HANDED_BACK
$ov.Range("C4").Value = $HANDED_BACK

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B3").Value = $HANDED_BACK
$zh.Range("B4").Value = $HANDED_BACK

# Latest Target File (E) / Latest Handback File (F) for row 3 (982bd13f...)
$zh.Range("E3").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.md"
$zh.Range("E3").Style = "Hyperlink"
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e9c090fb73bad951a50fe63180006ee3645059ab/e2e/982bd13f-222c-489b-8da0-a104e9cbbf60.md", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.md") | Out-Null

$zh.Range("F3").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf"
$zh.Range("F3").Style = "Hyperlink"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/792fa1781e7ad373b93cf685aeebc7b5dd82f97e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G) for row 3
$zh.Range("G3").Value = "2016-02-18 07:58:11"

# Row 4 (cc4d4c40...) inherits the same dependency info (it depends on 982bd13f...)
$zh.Range("E4").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.md"
$zh.Range("E4").Style = "Hyperlink"
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e9c090fb73bad951a50fe63180006ee3645059ab/e2e/982bd13f-222c-489b-8da0-a104e9cbbf60.md", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.md") | Out-Null

$zh.Range("F4").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf"
$zh.Range("F4").Style = "Hyperlink"
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/792fa1781e7ad373b93cf685aeebc7b5dd82f97e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.zh-cn.xlf") | Out-Null

$zh.Range("G4").Value = "2016-02-18 07:58:11"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B3").Value = $HANDED_BACK
$de.Range("B4").Value = $HANDED_BACK

$de.Range("E3").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.md"
$de.Range("E3").Style = "Hyperlink"
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e9c090fb73bad951a50fe63180006ee3645059ab/e2e/982bd13f-222c-489b-8da0-a104e9cbbf60.md", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.md") | Out-Null

$de.Range("F3").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf"
$de.Range("F3").Style = "Hyperlink"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fcbc632f0643712895668a1ff21e2bd10591e2fe/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf") | Out-Null

$de.Range("G3").Value = "2016-02-18 07:58:31"

$de.Range("E4").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.md"
$de.Range("E4").Style = "Hyperlink"
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e9c090fb73bad951a50fe63180006ee3645059ab/e2e/982bd13f-222c-489b-8da0-a104e9cbbf60.md", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.md") | Out-Null

$de.Range("F4").Value = "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf"
$de.Range("F4").Style = "Hyperlink"
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fcbc632f0643712895668a1ff21e2bd10591e2fe/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf", "", "", "982bd13f-222c-489b-8da0-a104e9cbbf60.8bf9e3e515168d35b3e4de2fc3acb986a37e4240.de-de.xlf") | Out-Null

$de.Range("G4").Value = "2016-02-18 07:58:31"

Write-Host "Handback report generated."
